# Apply commit changes to the "Avverkningsanmälningar" sheet:
# 1. Update column C (Förändrad) for rows 2-480 from 45182 to 45184.
# 2. Set row 480 height to 15 (customHeight) - matches new rows' formatting.
# 3. Append a new row (481) with the new cleared/updated entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed date) column for existing rows 2-480.
$ws.Range("C2:C480").Value = 45184

# 2. Ensure row 480 carries explicit row height metadata like the new row.
$ws.Rows.Item(480).RowHeight = 15

# 3. Add the new row 481 with its data.
$ws.Cells.Item(481, 1).Value = "A 43016-2023"            # A481 Beteckning
$ws.Cells.Item(481, 2).Value = 45182                      # B481 Datum
$ws.Cells.Item(481, 3).Value = 45184                      # C481 Förändrad
$ws.Cells.Item(481, 4).Value = "KRONOBERGS LÄN"           # D481 Län
$ws.Cells.Item(481, 5).Value = "ÄLMHULT"                  # E481 Kommun
$ws.Cells.Item(481, 7).Value = 0.8                        # G481 Area (ha)
$ws.Cells.Item(481, 8).Value = 0                          # H481 Fridlysta
$ws.Cells.Item(481, 9).Value = 0                          # I481 Signalarter
$ws.Cells.Item(481, 10).Value = 0                         # J481 NT
$ws.Cells.Item(481, 11).Value = 0                         # K481 VU
$ws.Cells.Item(481, 12).Value = 0                         # L481 EN
$ws.Cells.Item(481, 13).Value = 0                         # M481 CR
$ws.Cells.Item(481, 14).Value = 0                         # N481 RE
$ws.Cells.Item(481, 15).Value = 0                         # O481 Rödlistade
$ws.Cells.Item(481, 16).Value = 0                         # P481 Hotade
$ws.Cells.Item(481, 17).Value = 0                         # Q481 Alla arter

# Apply the same date number format (yyyy-mm-dd custom format, style used by column B/C).
$ws.Range("B481:C481").NumberFormat = "YYYY-MM-DD"

# R481 keeps the wrap-text style used throughout the rest of the column, but stays empty.
$ws.Cells.Item(481, 18).WrapText = $true
